$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("El Shaddai 105")

# Update quantity for the remaining line item (row 8, column C)
$ws.Range("C8").Value = 5

# Remove the two extra line-item rows (Dome Camera and Motorized Bullet entries)
$ws.Rows("9:10").Delete()

# Row 8 now only wraps to two lines instead of three, so shrink its height
$ws.Rows("8").RowHeight = 39.6

# Update the selection to match the final cursor position
$ws.Range("C9").Select()
